$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '62.803.29'
Set-TextValue $ws.Range("E2") '  -1.87%  '

Set-TextValue $ws.Range("D3") '3.056.51'
Set-TextValue $ws.Range("E3") '  -1.80%  '

Set-TextValue $ws.Range("E4") '  +0.24%  '

Set-TextValue $ws.Range("D5") '534.21'
Set-TextValue $ws.Range("E5") '  -4.50%  '

Set-TextValue $ws.Range("D6") '133.25'
Set-TextValue $ws.Range("E6") '  -4.08%  '

Set-TextValue $ws.Range("D7") '0.999'
Set-TextValue $ws.Range("E7") '  -0.07%  '

Set-TextValue $ws.Range("D8") '3.053.39'
Set-TextValue $ws.Range("E8") '  -1.75%  '

Set-TextValue $ws.Range("D9") '0.496'
Set-TextValue $ws.Range("E9") '  +0.79%  '

Set-TextValue $ws.Range("D10") '0.153'
Set-TextValue $ws.Range("E10") '  +1.06%  '

Set-TextValue $ws.Range("D11") '6.19'
Set-TextValue $ws.Range("E11") '  -8.84%  '

Set-TextValue $ws.Range("D12") '0.454'
Set-TextValue $ws.Range("E12") '  -0.79%  '

Set-TextValue $ws.Range("D13") '0.0000223'
Set-TextValue $ws.Range("E13") '  +2.59%  '

Set-TextValue $ws.Range("D14") '34.26'
Set-TextValue $ws.Range("E14") '  -4.10%  '

Set-TextValue $ws.Range("D15") '3.554.72'
Set-TextValue $ws.Range("E15") '  -1.66%  '

Set-TextValue $ws.Range("D16") '62.792.94'
Set-TextValue $ws.Range("E16") '  -1.87%  '

Set-TextValue $ws.Range("E17") '  -1.03%  '

Set-TextValue $ws.Range("D18") '3.068.10'
Set-TextValue $ws.Range("E18") '  -1.45%  '

Set-TextValue $ws.Range("D19") '6.61'
Set-TextValue $ws.Range("E19") '  -1.43%  '

Set-TextValue $ws.Range("D20") '481.80'
Set-TextValue $ws.Range("E20") '  -4.90%  '

Set-TextValue $ws.Range("D21") '13.25'
Set-TextValue $ws.Range("E21") '  -3.80%  '

Set-TextValue $ws.Range("D22") '0.694'
Set-TextValue $ws.Range("E22") '  -1.78%  '

Set-TextValue $ws.Range("D23") '7.13'
Set-TextValue $ws.Range("E23") '  -2.01%  '

Set-TextValue $ws.Range("D24") '79.24'
Set-TextValue $ws.Range("E24") '  +1.55%  '

Set-TextValue $ws.Range("D25") '12.11'
Set-TextValue $ws.Range("E25") '  -3.05%  '

Set-TextValue $ws.Range("D26") '0.997'
Set-TextValue $ws.Range("E26") '  -0.25%  '

Set-TextValue $ws.Range("D27") '2.68'
Set-TextValue $ws.Range("E27") '  -3.92%  '

Set-TextValue $ws.Range("D28") '8.09'
Set-TextValue $ws.Range("E28") '  -4.76%  '

Set-TextValue $ws.Range("D29") '1.00'
Set-TextValue $ws.Range("E29") '  +0.36%  '

Set-TextValue $ws.Range("D30") '25.86'
Set-TextValue $ws.Range("E30") '  -1.69%  '

Set-TextValue $ws.Range("D31") '1.86'
Set-TextValue $ws.Range("E31") '  -9.62%  '

Set-TextValue $ws.Range("E32") '  -1.02%  '

Set-TextValue $ws.Range("D33") '2.37'
Set-TextValue $ws.Range("E33") '  -7.61%  '

Set-TextValue $ws.Range("D34") '56.37'
Set-TextValue $ws.Range("E34") '  +2.21%  '

Set-TextValue $ws.Range("D35") '5.34'
Set-TextValue $ws.Range("E35") '  +1.83%  '

Set-TextValue $ws.Range("D36") '5.93'
Set-TextValue $ws.Range("E36") '  -0.10%  '

Set-TextValue $ws.Range("D37") '473.60'
Set-TextValue $ws.Range("E37") '  -13.58%  '

Set-TextValue $ws.Range("B38") 'Maker'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D38") '3.093.41'
Set-TextValue $ws.Range("E38") '  +0.59%  '

Set-TextValue $ws.Range("B39") 'VeChain'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D39") '0.0393'
Set-TextValue $ws.Range("E39") '  -5.85%  '

Set-TextValue $ws.Range("D40") '0.0793'
Set-TextValue $ws.Range("E40") '  -1.31%  '

Set-TextValue $ws.Range("E41") '  -3.56%  '

Set-TextValue $ws.Range("D42") '8.07'
Set-TextValue $ws.Range("E42") '  -1.05%  '

Set-TextValue $ws.Range("D43") '2.60'
Set-TextValue $ws.Range("E43") '  -1.29%  '

Set-TextValue $ws.Range("D44") '0.252'
Set-TextValue $ws.Range("E44") '  -1.98%  '

Set-TextValue $ws.Range("E45") '  +0.06%  '

Set-TextValue $ws.Range("D46") '0.0₃0542'
Set-TextValue $ws.Range("E46") '  +7.30%  '

Set-TextValue $ws.Range("D47") '2.03'
Set-TextValue $ws.Range("E47") '  -4.63%  '

Set-TextValue $ws.Range("B48") 'Monero'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D48") '120.76'
Set-TextValue $ws.Range("E48") '  -0.38%  '

Set-TextValue $ws.Range("B49") 'InjectiveProtocol'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D49") '24.52'
Set-TextValue $ws.Range("E49") '  +0.41%  '

Set-TextValue $ws.Range("E50") '  +0.45%  '

Set-TextValue $ws.Range("D51") '2.33'
Set-TextValue $ws.Range("E51") '  +5.15%  '

